$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 72.266001
$ws.Range("H2").Value2 = 216.798003
$ws.Range("I2").Value2 = 0.2949652269937106
$ws.Range("J2").Value2 = 0.2949652269937106
$ws.Range("M2").Value2 = 568.5612486666666
$ws.Range("N2").Value2 = 1705.683746
$ws.Range("O2").Value2 = 0.6737621253161296
$ws.Range("P2").Value2 = 0.6737621253161296
$ws.Range("Q2").Value2 = 41087.64776470658
$ws.Range("R2").Value2 = 369788.8298823592
$ws.Range("S2").Value2 = 0.198736398233637
$ws.Range("T2").Value2 = 0.198736398233637

$ws.Range("G3").Value2 = 72.266001
$ws.Range("H3").Value2 = 216.798003
$ws.Range("I3").Value2 = 0.2949652269937106
$ws.Range("J3").Value2 = 0.2949652269937106
$ws.Range("M3").Value2 = 88.00803400000001
$ws.Range("N3").Value2 = 264.024102
$ws.Range("O3").Value2 = 0.1042921587987053
$ws.Range("P3").Value2 = 0.1042921587987053
$ws.Range("Q3").Value2 = 6359.988673052035
$ws.Range("R3").Value2 = 57239.89805746831
$ws.Range("S3").Value2 = 0.03076256029372423
$ws.Range("T3").Value2 = 0.03076256029372423

$ws.Range("G4").Value2 = 72.266001
$ws.Range("H4").Value2 = 216.798003
$ws.Range("I4").Value2 = 0.2949652269937106
$ws.Range("J4").Value2 = 0.2949652269937106
$ws.Range("M4").Value2 = 187.2912243333334
$ws.Range("N4").Value2 = 561.873673
$ws.Range("O4").Value2 = 0.2219457158851651
$ws.Range("P4").Value2 = 0.2219457158851651
$ws.Range("Q4").Value2 = 13534.78780496389
$ws.Range("R4").Value2 = 121813.090244675
$ws.Range("S4").Value2 = 0.06546626846634931
$ws.Range("T4").Value2 = 0.06546626846634931

$ws.Range("H5").Value2 = 410.023338
$ws.Range("I5").Value2 = 0.5578585839920717
$ws.Range("J5").Value2 = 0.5578585839920718
$ws.Range("M5").Value2 = 568.5612486666666
$ws.Range("N5").Value2 = 1705.683746
$ws.Range("O5").Value2 = 0.6737621253161296
$ws.Range("P5").Value2 = 0.6737621253161296
$ws.Range("Q5").Value2 = 77707.7936785849
$ws.Range("R5").Value2 = 699370.1431072641
$ws.Range("S5").Value2 = 0.3758639851763448
$ws.Range("T5").Value2 = 0.3758639851763448

$ws.Range("H6").Value2 = 410.023338
$ws.Range("I6").Value2 = 0.5578585839920717
$ws.Range("J6").Value2 = 0.5578585839920718
$ws.Range("M6").Value2 = 88.00803400000001
$ws.Range("O6").Value2 = 0.1042921587987053
$ws.Range("P6").Value2 = 0.1042921587987053
$ws.Range("Q6").Value2 = 12028.44929049916
$ws.Range("S6").Value2 = 0.05818027602892204
$ws.Range("T6").Value2 = 0.05818027602892206

$ws.Range("H7").Value2 = 410.023338
$ws.Range("I7").Value2 = 0.5578585839920717
$ws.Range("J7").Value2 = 0.5578585839920718
$ws.Range("M7").Value2 = 187.2912243333334
$ws.Range("N7").Value2 = 561.873673
$ws.Range("O7").Value2 = 0.2219457158851651
$ws.Range("P7").Value2 = 0.2219457158851651
$ws.Range("Q7").Value2 = 25597.92432642006
$ws.Range("R7").Value2 = 230381.3189377805
$ws.Range("S7").Value2 = 0.1238143227868048
$ws.Range("T7").Value2 = 0.1238143227868049

$ws.Range("G8").Value2 = 36.057927
$ws.Range("H8").Value2 = 108.173781
$ws.Range("I8").Value2 = 0.1471761890142177
$ws.Range("J8").Value2 = 0.1471761890142177
$ws.Range("M8").Value2 = 568.5612486666666
$ws.Range("N8").Value2 = 1705.683746
$ws.Range("O8").Value2 = 0.6737621253161296
$ws.Range("P8").Value2 = 0.6737621253161296
$ws.Range("Q8").Value2 = 20501.13999945151
$ws.Range("R8").Value2 = 184510.2599950636
$ws.Range("S8").Value2 = 0.09916174190614771
$ws.Range("T8").Value2 = 0.09916174190614771

$ws.Range("G9").Value2 = 36.057927
$ws.Range("H9").Value2 = 108.173781
$ws.Range("I9").Value2 = 0.1471761890142177
$ws.Range("J9").Value2 = 0.1471761890142177
$ws.Range("M9").Value2 = 88.00803400000001
$ws.Range("O9").Value2 = 0.1042921587987053
$ws.Range("P9").Value2 = 0.1042921587987053
$ws.Range("Q9").Value2 = 3173.387265385518
$ws.Range("R9").Value2 = 28560.48538846966
$ws.Range("S9").Value2 = 0.01534932247605907
$ws.Range("T9").Value2 = 0.01534932247605907

$ws.Range("G10").Value2 = 36.057927
$ws.Range("H10").Value2 = 108.173781
$ws.Range("I10").Value2 = 0.1471761890142177
$ws.Range("J10").Value2 = 0.1471761890142177
$ws.Range("M10").Value2 = 187.2912243333334
$ws.Range("N10").Value2 = 561.873673
$ws.Range("O10").Value2 = 0.2219457158851651
$ws.Range("P10").Value2 = 0.2219457158851651
$ws.Range("Q10").Value2 = 6753.333294751958
$ws.Range("R10").Value2 = 60779.99965276761
$ws.Range("S10").Value2 = 0.03266512463201092
$ws.Range("T10").Value2 = 0.03266512463201091

